$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTX")

# Row 17: Gross Margin
$ws.Range("D17").Value = 0.1804
$ws.Range("E17").Value = 0.2128
$ws.Range("F17").Value = 0.2417
$ws.Range("G17").Value = 0.2371

# Row 18: EBIT Margin
$ws.Range("D18").Value = -0.0385
$ws.Range("E18").Value = -0.0205
$ws.Range("F18").Value = 0.0921
$ws.Range("G18").Value = 0.1084

# Row 19: EBT margin
$ws.Range("D19").Value = -0.0542
$ws.Range("E19").Value = -0.0377
$ws.Range("F19").Value = 0.0726
$ws.Range("G19").Value = 0.0916

# Row 20: Net Profit Margin
$ws.Range("D20").Value = -0.0489
$ws.Range("E20").Value = -0.0339
$ws.Range("F20").Value = 0.0908
$ws.Range("G20").Value = 0.1221

# Row 21: Free Cash Flow Margin
$ws.Range("D21").Value = 0.0282
$ws.Range("E21").Value = 0.0307
$ws.Range("F21").Value = 0.0621
$ws.Range("G21").Value = 0.0794

# Row 32: EBITDA Margin
$ws.Range("D32").Value = 0.1079
$ws.Range("E32").Value = 0.1354
$ws.Range("F32").Value = 0.1669
$ws.Range("G32").Value = 0.1681

# Row 33: Operating Cash Flow Margin
$ws.Range("D33").Value = 0.07
$ws.Range("E33").Value = 0.0778
$ws.Range("F33").Value = 0.1102
$ws.Range("G33").Value = 0.1284
